$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Aarau -> Chiasso, Road)
$ws.Range("E2").Value = 3.287119813793048
$ws.Range("F2").Value = 255.74848
$ws.Range("G2").Value = 2.939999999999987
$ws.Range("H2").Value = 1

# Row 3 (Aarau -> Chiasso, Intermodal)
$ws.Range("G3").Value = 144.06

# Row 4 (Chiasso -> Aarau, Road)
$ws.Range("E4").Value = 112
$ws.Range("F4").Value = 249.65772
$ws.Range("G4").Value = 112
$ws.Range("H4").Value = 1

# Row 5 (Chiasso -> Aarau, Intermodal)
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0

# Row 6 (Aarau -> Stabio, Road)
$ws.Range("D6").Value = 0
$ws.Range("H6").Value = 1

# Row 7 (Aarau -> Stabio, Intermodal)
$ws.Range("D7").Value = 0

# Row 8 (Stabio -> Aarau, Road)
$ws.Range("D8").Value = 0
$ws.Range("H8").Value = 1

# Row 9 (Stabio -> Aarau, Intermodal)
$ws.Range("D9").Value = 0

# Row 10 (Aarau -> Visp, Road)
$ws.Range("E10").Value = 63
$ws.Range("F10").Value = 229.86704
$ws.Range("G10").Value = 63
$ws.Range("H10").Value = 1

# Row 11 (Aarau -> Visp, Intermodal)
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0

# Row 12 (Visp -> Aarau, Road)
$ws.Range("F12").Value = 250.64208
$ws.Range("H12").Value = 1

# Row 14 (Baselwolf -> Chiasso, Road)
$ws.Range("E14").Value = 1764
$ws.Range("F14").Value = 334.4031187301584
$ws.Range("G14").Value = 1764
$ws.Range("H14").Value = 1

# Row 15 (Baselwolf -> Chiasso, Intermodal)
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0

# Row 16 (Chiasso -> Baselwolf, Road)
$ws.Range("E16").Value = 1442
$ws.Range("F16").Value = 332.5168877669899
$ws.Range("G16").Value = 1442
$ws.Range("H16").Value = 1

# Row 17 (Chiasso -> Baselwolf, Intermodal)
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
